# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 2 (top of the data block),
# pushing the existing rows 2-21 down to rows 3-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 21
$lastCol = 20   # columns A..T

# Shift existing data rows down by one, starting from the bottom so we
# never overwrite a row before it has been copied.
for ($r = $lastRow; $r -ge 2; $r--) {
    $dstRow = $r + 1
    for ($col = 1; $col -le $lastCol; $col++) {
        $srcCell = $ws.Cells.Item($r, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)
        $dstCell.Value = $srcCell.Value2
    }
}

# Row 22 is brand new (previously unused), so it needs the "Fecha" date
# number format applied to column D to match the rest of the column.
$ws.Cells.Item($lastRow + 1, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Populate the new row 2 with the new weekly record.
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(2, 3).Value = 'Maule'
$ws.Cells.Item(2, 4).Value = 45083
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 'Fruta'
$ws.Cells.Item(2, 7).Value = 100104
$ws.Cells.Item(2, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(2, 9).Value = 100104001
$ws.Cells.Item(2, 10).Value = 'Granada'
$ws.Cells.Item(2, 11).Value = 'Wonderfull'
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 120
$ws.Cells.Item(2, 14).Value = 17000
$ws.Cells.Item(2, 15).Value = 17000
$ws.Cells.Item(2, 16).Value = 17000
$ws.Cells.Item(2, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 944
$ws.Cells.Item(2, 20).Value = 18
